$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) date values from 45186 (2023-09-17) to 45188 (2023-09-19)
# for all data rows (2 through 176), preserving existing cell formatting.
for ($r = 2; $r -le 176; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45186) {
        $cell.Value2 = 45188
    }
}
